$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.933.92"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.903.92"
$ws.Range("E3").Value = "  -3.46%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.17"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4591"
$ws.Range("E7").Value = "  -1.49%  "
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07708"
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9771"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("E11").Value = "  -3.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.926.44"
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.927"
$ws.Range("E13").Value = "  -3.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.653"
$ws.Range("E14").Value = "  -3.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07054"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "83.60"
$ws.Range("E17").Value = "  -4.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009445"
$ws.Range("E18").Value = "  -5.01%  "
$ws.Range("E19").Value = "  -3.76%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.940.16"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.295"
$ws.Range("E22").Value = "  -4.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("E23").Value = "  -2.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.157.91"
$ws.Range("E24").Value = "  -2.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.095"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.16"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.06"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.648"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.30"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.837"
$ws.Range("E30").Value = "  -3.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09266"
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("E32").Value = "  -3.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.080"
$ws.Range("E33").Value = "  -2.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.243"
$ws.Range("E34").Value = "  -5.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.949"
$ws.Range("E35").Value = "  -6.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05713"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.146"
$ws.Range("E37").Value = "  -2.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.002"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02039"
$ws.Range("E39").Value = "  -3.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5488"
$ws.Range("E40").Value = "  -4.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.386"
$ws.Range("E41").Value = "  -4.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1754"
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.282"
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.773"
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5165"
$ws.Range("E45").Value = "  -3.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.24"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06808"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("E48").Value = "  -6.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000002568"
$ws.Range("E49").Value = "  -17.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.22"
$ws.Range("E50").Value = "  -3.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.769"
$ws.Range("E51").Value = "  -3.22%  "
